$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 223
$ws1.Range("F3").Value = 161
$ws1.Range("F4").Value = 0

# Sheet "全部类型" (All types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 223
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 147
$ws4.Range("F5").Value = 0
